# Feedback details. Remove stimulus name from item page.
#
# Updates the FEEDBACK_SHORT row (row 17) on Sheet1 with new, longer
# feedback text for the English (col B) and German (col C) locales, and
# adjusts the row height to fit the new (wrapped) text. Also updates the
# last-saved cell selection to reflect where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FEEDBACK_SHORT strings (replace the old, shorter "You rated ..." /
# "Sie haben ..." texts referenced by row 17).
$ws.Range("B17").Value = "You finished the Expressivity Comparison Test with {{num_correct}} out of {{num_questions}} ({{perc_correct}}%) correct answers. "
$ws.Range("C17").Value = "Sie haben den Test zum Ausdruck in der Musik mit {{num_correct}} von {{num_questions}} ({{perc_correct}}%) richtigen Antworten beendet."

# The longer text now wraps across more lines, so the row grows from 30 to 45.
$ws.Rows.Item(17).RowHeight = 45

# Reflect the final cursor position/selection as last saved.
$ws.Range("C18").Select()
